$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '66.712.96'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  +2.59%  '
$ws.Range('E2').Style = 'Normal'
$ws.Range('D3').Value = '3.731.51'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  +6.57%  '
$ws.Range('E3').Style = 'Normal'
$ws.Range('E4').Value = '  +0.31%  '
$ws.Range('E4').Style = 'Normal'
$ws.Range('D5').Value = '''419.34'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +0.43%  '
$ws.Range('E5').Style = 'Normal'
$ws.Range('D6').Value = '''131.57'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -0.21%  '
$ws.Range('E6').Style = 'Normal'
$ws.Range('D7').Value = '3.719.63'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +6.53%  '
$ws.Range('E7').Style = 'Normal'
$ws.Range('D8').Value = '''0.650'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -1.87%  '
$ws.Range('E8').Style = 'Normal'
$ws.Range('E9').Value = '  +0.13%  '
$ws.Range('E9').Style = 'Normal'
$ws.Range('D10').Value = '''0.770'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -1.32%  '
$ws.Range('E10').Style = 'Normal'
$ws.Range('D11').Value = '''0.182'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +12.42%  '
$ws.Range('E11').Style = 'Normal'
$ws.Range('D12').Value = '''0.0000404'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +56.39%  '
$ws.Range('E12').Style = 'Normal'
$ws.Range('D13').Value = '''42.86'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -0.96%  '
$ws.Range('E13').Style = 'Normal'
$ws.Range('D14').Value = '''10.51'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +5.84%  '
$ws.Range('E14').Style = 'Normal'
$ws.Range('D15').Value = '4.321.28'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +6.72%  '
$ws.Range('E15').Style = 'Normal'
$ws.Range('E16').Value = '  -0.93%  '
$ws.Range('E16').Style = 'Normal'
$ws.Range('D17').Value = '''20.77'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +1.86%  '
$ws.Range('E17').Style = 'Normal'
$ws.Range('D18').Value = '3.716.55'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +6.28%  '
$ws.Range('E18').Style = 'Normal'
$ws.Range('D19').Value = '''13.27'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +5.34%  '
$ws.Range('E19').Style = 'Normal'
$ws.Range('E20').Value = '  +4.43%  '
$ws.Range('E20').Style = 'Normal'
$ws.Range('D21').Value = '66.813.03'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +3.00%  '
$ws.Range('E21').Style = 'Normal'
$ws.Range('D22').Value = '''444.83'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -2.25%  '
$ws.Range('E22').Style = 'Normal'
$ws.Range('E23').Value = '  +24.57%  '
$ws.Range('E23').Style = 'Normal'
$ws.Range('D24').Value = '''89.81'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -0.36%  '
$ws.Range('E24').Style = 'Normal'
$ws.Range('E25').Value = '  -0.95%  '
$ws.Range('E25').Style = 'Normal'
$ws.Range('D26').Value = '''38.53'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +13.41%  '
$ws.Range('E26').Style = 'Normal'
$ws.Range('D27').Value = '''10.22'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +3.02%  '
$ws.Range('E27').Style = 'Normal'
$ws.Range('D28').Value = '''3.34'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -0.72%  '
$ws.Range('E28').Style = 'Normal'
$ws.Range('E29').Value = '  +4.40%  '
$ws.Range('E29').Style = 'Normal'
$ws.Range('D30').Value = '''12.78'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +2.40%  '
$ws.Range('E30').Style = 'Normal'
$ws.Range('E31').Value = '  +9.68%  '
$ws.Range('E31').Style = 'Normal'
$ws.Range('E32').Value = '  +1.98%  '
$ws.Range('E32').Style = 'Normal'
$ws.Range('D33').Value = '''7.25'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -2.89%  '
$ws.Range('E33').Style = 'Normal'
$ws.Range('E34').Value = '  +0.80%  '
$ws.Range('E34').Style = 'Normal'
$ws.Range('D35').Value = '''41.99'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +6.06%  '
$ws.Range('E35').Style = 'Normal'
$ws.Range('D36').Value = '''57.02'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -0.62%  '
$ws.Range('E36').Style = 'Normal'
$ws.Range('E37').Value = '  +0.09%  '
$ws.Range('E37').Style = 'Normal'
$ws.Range('D38').Value = '''0.0495'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -2.12%  '
$ws.Range('E38').Style = 'Normal'
$ws.Range('D39').Value = '0.0₃0750'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +7.40%  '
$ws.Range('E39').Style = 'Normal'
$ws.Range('D40').Value = '''0.151'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -3.25%  '
$ws.Range('E40').Style = 'Normal'
$ws.Range('D41').Value = '''3.05'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +30.89%  '
$ws.Range('E41').Style = 'Normal'
$ws.Range('D42').Value = '''29.17'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +35.59%  '
$ws.Range('E42').Style = 'Normal'
$ws.Range('D43').Value = '''0.998'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +0.10%  '
$ws.Range('E43').Style = 'Normal'
$ws.Range('E44').Value = '  +4.55%  '
$ws.Range('E44').Style = 'Normal'
$ws.Range('D45').Value = '''3.24'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +32.07%  '
$ws.Range('E45').Style = 'Normal'
$ws.Range('E46').Value = '  +6.33%  '
$ws.Range('E46').Style = 'Normal'
$ws.Range('D47').Value = '''146.41'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -0.92%  '
$ws.Range('E47').Style = 'Normal'
$ws.Range('E48').Value = '  -3.96%  '
$ws.Range('E48').Style = 'Normal'
$ws.Range('B49').Value = 'Stacks'
$ws.Range('B49').Style = 'Normal'
$ws.Range('C49').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('C49').Style = 'Normal'
$ws.Range('D49').Value = '''2.89'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -5.55%  '
$ws.Range('E49').Style = 'Normal'
$ws.Range('B50').Value = 'NEARProtocol'
$ws.Range('B50').Style = 'Normal'
$ws.Range('C50').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('C50').Style = 'Normal'
$ws.Range('D50').Value = '''4.37'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -3.51%  '
$ws.Range('E50').Style = 'Normal'
$ws.Range('D51').Value = '''0.307'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -2.29%  '
$ws.Range('E51').Style = 'Normal'
